$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text runs): Volume/Number and week-of dates ---
# A8: "...Number  23" -> "...Number  24" (chars 21-22 of the merged text)
$ws.Range("A8").Characters(21, 2).Text = "24"
# C9: "Report Covering the Week  6/2/2025  Through  6/8/2025"
#     -> "Report Covering the Week  6/9/2025  Through  6/15/2025"
$ws.Range("C9").Characters(27, 6).Text = "6/9/2025"
$ws.Range("C9").Characters(45, 7).Text = "6/15/2025"

# --- Cells switching between "n/a" placeholder (shared text) and numeric value ---
# text -> number (style must become the numeric style of the column)
$ws.Range("D19").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C17").Value = 4
$ws.Range("D19").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = 4
$ws.Range("E19").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E18").Value = -50
$ws.Range("D19").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value = 1
$ws.Range("E19").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").Value = -100

# number -> "n/a" placeholder text (style must become the text/General style)
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial(-4122)

# --- Plain numeric updates (weekly crime-stat figures refresh) ---
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 35
$ws.Range("J16").Value = 36
$ws.Range("K16").Value = -2.777777777777
$ws.Range("L16").Value = 6.060606060606
$ws.Range("M16").Value = -20.454545454545
$ws.Range("N16").Value = -89.096573208722
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 33.333333333333
$ws.Range("F17").Value = 11
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = 22.222222222222
$ws.Range("I17").Value = 72
$ws.Range("J17").Value = 45
$ws.Range("K17").Value = 60
$ws.Range("L17").Value = 56.521739130434
$ws.Range("M17").Value = 157.142857142857
$ws.Range("N17").Value = 28.571428571428
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = 50
$ws.Range("I18").Value = 72
$ws.Range("J18").Value = 45
$ws.Range("K18").Value = 60
$ws.Range("L18").Value = 30.909090909090
$ws.Range("M18").Value = 26.315789473684
$ws.Range("N18").Value = -87.919463087248
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -30
$ws.Range("F19").Value = 23
$ws.Range("G19").Value = 38
$ws.Range("H19").Value = -39.473684210526
$ws.Range("I19").Value = 165
$ws.Range("J19").Value = 202
$ws.Range("K19").Value = -18.316831683168
$ws.Range("L19").Value = -25.339366515837
$ws.Range("M19").Value = 4.430379746835
$ws.Range("N19").Value = -62.669683257918
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 12
$ws.Range("H20").Value = -53.846153846153
$ws.Range("I20").Value = 72
$ws.Range("J20").Value = 78
$ws.Range("K20").Value = -7.692307692307
$ws.Range("L20").Value = 30.909090909090
$ws.Range("M20").Value = 44
$ws.Range("N20").Value = -95.411089866156
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = -22.727272727272
$ws.Range("F21").Value = 60
$ws.Range("G21").Value = 85
$ws.Range("H21").Value = -29.411764705882
$ws.Range("I21").Value = 422
$ws.Range("J21").Value = 411
$ws.Range("K21").Value = 2.676399026763
$ws.Range("L21").Value = 2.427184466019
$ws.Range("M21").Value = 24.852071005917
$ws.Range("N21").Value = -85.876840696117
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = -100
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = -80
$ws.Range("J22").Value = 20
$ws.Range("K22").Value = -25
$ws.Range("C24").Value = 50
$ws.Range("D24").Value = 41
$ws.Range("E24").Value = 21.951219512195
$ws.Range("F24").Value = 166
$ws.Range("G24").Value = 136
$ws.Range("H24").Value = 22.058823529411
$ws.Range("I24").Value = 936
$ws.Range("J24").Value = 736
$ws.Range("K24").Value = 27.173913043478
$ws.Range("L24").Value = 29.281767955801
$ws.Range("M24").Value = 119.718309859155
$ws.Range("C25").Value = 39
$ws.Range("D25").Value = 27
$ws.Range("E25").Value = 44.444444444444
$ws.Range("F25").Value = 133
$ws.Range("G25").Value = 87
$ws.Range("H25").Value = 52.873563218390
$ws.Range("I25").Value = 736
$ws.Range("J25").Value = 535
$ws.Range("K25").Value = 37.570093457943
$ws.Range("L25").Value = 40.726577437858
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = -60
$ws.Range("I26").Value = 140
$ws.Range("J26").Value = 129
$ws.Range("K26").Value = 8.527131782945
$ws.Range("L26").Value = 26.126126126126
$ws.Range("M26").Value = 40
$ws.Range("G27").Value = 2
$ws.Range("J27").Value = 10
$ws.Range("K27").Value = -10
$ws.Range("E28").Value = -100
$ws.Range("J28").Value = 15
$ws.Range("K28").Value = -6.666666666666
$ws.Range("L28").Value = 0

$excel.CutCopyMode = 0
